$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.098.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.474.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.37%  '

$ws.Range("E4").Value = '  -0.82%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.59%  '

$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.539'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.50%  '

$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.473.14'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.85%  '

$ws.Range("E10").Value = '  +0.58%  '

$ws.Range("E11").Value = '  +1.48%  '

$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.92%  '

$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.91%  '

$ws.Range("E15").Value = '  +0.85%  '

$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.923.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.29%  '

$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.048.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.32%  '

$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.465.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("E19").Value = '  +3.71%  '

$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.53%  '

$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '329.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.56%  '

$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("E25").Value = '  +0.93%  '

$ws.Range("D26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '668.55'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.65%  '

$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.18%  '

$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0993'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.60%  '

$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.593.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.09%  '

$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.46%  '

$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.57%  '

$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.26%  '

$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.60%  '

$ws.Range("E34").Value = '  -4.00%  '

$ws.Range("E35").Value = '  +3.35%  '

$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.61%  '

$ws.Range("E38").Value = '  +0.96%  '

$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.372'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.55%  '

$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.80'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '151.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.77%  '

$ws.Range("E42").Value = '  -0.45%  '

$ws.Range("D43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.61%  '

$ws.Range("E45").Value = '  +6.89%  '

$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '153.48'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.65%  '

$ws.Range("E47").Value = '  +19.23%  '

$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.64%  '

$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.606'
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = '  -0.95%  '
